{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Effective changes reproduced from the diff:\n//  1) The Word-managed \"_GoBack\" bookmark (marks the last edit location) is\n//     removed from its old spot (the empty paragraph right after the author\n//     list) and re-inserted inside the title run, right after\n//     \"Teasing apart mou\" (splitting that run in two, exactly like the diff).\n//  2) \"Manual to use the R script provided in Dryad\" becomes\n//     \"Manual to use the R script provided in GitHub\", with \"GitHub\" in its\n//     own run wrapped in spell-check proofErr markers (as in the diff).\n\nconst body = context.document.body;\n\n// --- 1) Relocate the \"_GoBack\" bookmark -----------------------------------\n// Remove the existing (hidden) \"_GoBack\" bookmark, wherever it currently is.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Find the split point inside the title and drop the bookmark there.\nconst titleHits = body.search(\"Teasing apart mou\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  const splitPoint = titleHits.items[0].getRange(\"End\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 2) Replace \"Dryad\" with \"GitHub\" --------------------------------------\nconst dryadHits = body.search(\"Dryad\", { matchCase: true });\ndryadHits.load(\"items\");\nawait context.sync();\n\nif (dryadHits.items.length > 0) {\n  // Keep the same run formatting (bold, blue, size 28, Times New Roman) and\n  // wrap \"GitHub\" with spellStart/spellEnd proofErr markers, same as the\n  // word it replaces would have produced through Word's own spell check.\n  const ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:b/><w:color w:val=\"3366FF\"/><w:sz w:val=\"28\"/></w:rPr><w:t>GitHub</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n  dryadHits.items[0].insertOoxml(ooxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Effective changes reproduced from the diff:\n#  1) The Word-managed \"_GoBack\" bookmark (marks the last edit location) is\n#     removed from its old spot (the empty paragraph right after the author\n#     list) and re-inserted inside the title run, right after\n#     \"Teasing apart mou\" (splitting that run in two, exactly like the diff).\n#  2) \"Manual to use the R script provided in Dryad\" -> \"...provided in GitHub\".\n\n$d = $word.ActiveDocument\n\n# --- 1) Relocate the \"_GoBack\" bookmark ------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$titleRange = $d.Content\n$titleRange.Find.Execute(\"Teasing apart mou\") | Out-Null\n$titleRange.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd)\n$d.Bookmarks.Add(\"_GoBack\", $titleRange) | Out-Null\n\n# --- 2) Replace \"Dryad\" with \"GitHub\" ---------------------------------------\n$dryadRange = $d.Content\nif ($dryadRange.Find.Execute(\"Dryad\")) {\n    $dryadRange.Text = \"GitHub\"\n}\n"}
